$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (values changed per diff) ---
# Row 4
$ws.Cells.Item(4, 3).Value = -2
$ws.Cells.Item(4, 4).Value = 45829.29943059995
$ws.Cells.Item(4, 5).Value = -2
$ws.Cells.Item(4, 6).Value = 45828.62064814815

# Row 12
$ws.Cells.Item(12, 3).Value = 64
$ws.Cells.Item(12, 4).Value = 45829.29943059795
$ws.Cells.Item(12, 5).Value = 64
$ws.Cells.Item(12, 6).Value = 45828.52197916667

# Row 24
$ws.Cells.Item(24, 3).Value = 49
$ws.Cells.Item(24, 4).Value = 45829.29943059705
$ws.Cells.Item(24, 5).Value = 49
$ws.Cells.Item(24, 6).Value = 45828.52178240741

# Row 33
$ws.Cells.Item(33, 3).Value = 3001
$ws.Cells.Item(33, 4).Value = 45829.29943060142
$ws.Cells.Item(33, 5).Value = 3001
$ws.Cells.Item(33, 6).Value = 45828.63898148148

# Row 61
$ws.Cells.Item(61, 3).Value = 4
$ws.Cells.Item(61, 4).Value = 45831.31672868288
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = 45829.58444444444

# Row 82
$ws.Cells.Item(82, 3).Value = 63
$ws.Cells.Item(82, 4).Value = 45829.29943059525
$ws.Cells.Item(82, 5).Value = 63
$ws.Cells.Item(82, 6).Value = 45828.4925

# Row 126
$ws.Cells.Item(126, 3).Value = 13
$ws.Cells.Item(126, 4).Value = 45831.31672868304
$ws.Cells.Item(126, 5).Value = 13
$ws.Cells.Item(126, 6).Value = 45829.58472222222
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = "Consistente"

# Row 151
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 4).Value = 45829.29943059936
$ws.Cells.Item(151, 5).Value = 2
$ws.Cells.Item(151, 6).Value = 45828.54791666667

# Row 161
$ws.Cells.Item(161, 3).Value = 77
$ws.Cells.Item(161, 4).Value = 45829.29943060024
$ws.Cells.Item(161, 5).Value = 77
$ws.Cells.Item(161, 6).Value = 45828.62064814815

# Row 195
$ws.Cells.Item(195, 3).Value = -3
$ws.Cells.Item(195, 4).Value = 45829.29943060406
$ws.Cells.Item(195, 5).Value = -3
$ws.Cells.Item(195, 6).Value = 45828.6547337963

# Row 206
$ws.Cells.Item(206, 3).Value = 8
$ws.Cells.Item(206, 4).Value = 45829.29943059968
$ws.Cells.Item(206, 5).Value = 8
$ws.Cells.Item(206, 6).Value = 45828.54859953704

# Row 263
$ws.Cells.Item(263, 3).Value = 21
$ws.Cells.Item(263, 4).Value = 45831.31672867625
$ws.Cells.Item(263, 5).Value = 21
$ws.Cells.Item(263, 6).Value = 45829.47876157407

# Row 287
$ws.Cells.Item(287, 3).Value = 487
$ws.Cells.Item(287, 4).Value = 45829.29943060171
$ws.Cells.Item(287, 5).Value = 487
$ws.Cells.Item(287, 6).Value = 45828.63898148148

# Row 300
$ws.Cells.Item(300, 3).Value = -2
$ws.Cells.Item(300, 4).Value = 45829.29943059734
$ws.Cells.Item(300, 5).Value = -2
$ws.Cells.Item(300, 6).Value = 45828.52178240741

# Row 379
$ws.Cells.Item(379, 3).Value = 17
$ws.Cells.Item(379, 4).Value = 45829.29943060053
$ws.Cells.Item(379, 5).Value = 17
$ws.Cells.Item(379, 6).Value = 45828.62064814815

# Row 404
$ws.Cells.Item(404, 3).Value = 195
$ws.Cells.Item(404, 4).Value = 45829.299430602
$ws.Cells.Item(404, 5).Value = 195
$ws.Cells.Item(404, 6).Value = 45828.63898148148

# Row 441
$ws.Cells.Item(441, 3).Value = 42
$ws.Cells.Item(441, 4).Value = 45829.29943059586
$ws.Cells.Item(441, 5).Value = 42
$ws.Cells.Item(441, 6).Value = 45828.52141203704

# Row 461
$ws.Cells.Item(461, 3).Value = 150
$ws.Cells.Item(461, 4).Value = 45829.29943060082
$ws.Cells.Item(461, 5).Value = 150
$ws.Cells.Item(461, 6).Value = 45828.62064814815

# Row 480
$ws.Cells.Item(480, 3).Value = 212
$ws.Cells.Item(480, 4).Value = 45829.2994306023
$ws.Cells.Item(480, 5).Value = 212
$ws.Cells.Item(480, 6).Value = 45828.63898148148

# Row 507
$ws.Cells.Item(507, 3).Value = 67
$ws.Cells.Item(507, 4).Value = 45831.3167286812
$ws.Cells.Item(507, 5).Value = 67
$ws.Cells.Item(507, 6).Value = 45829.55359953704

# Row 527
$ws.Cells.Item(527, 3).Value = 17
$ws.Cells.Item(527, 4).Value = 45829.29943059615
$ws.Cells.Item(527, 5).Value = 17
$ws.Cells.Item(527, 6).Value = 45828.52141203704

# Row 533
$ws.Cells.Item(533, 4).Value = 45831.31672868165

# Row 534
$ws.Cells.Item(534, 3).Value = 1233
$ws.Cells.Item(534, 4).Value = 45831.31672868149
$ws.Cells.Item(534, 5).Value = 1233
$ws.Cells.Item(534, 6).Value = 45829.56215277778

# Row 535
$ws.Cells.Item(535, 4).Value = 45831.31672868136

# Row 603
$ws.Cells.Item(603, 3).Value = 4
$ws.Cells.Item(603, 4).Value = 45829.29943060602
$ws.Cells.Item(603, 5).Value = 4
$ws.Cells.Item(603, 6).Value = 45828.71103009259

# Row 685
$ws.Cells.Item(685, 3).Value = 4
$ws.Cells.Item(685, 4).Value = 45829.29943060628
$ws.Cells.Item(685, 5).Value = 4
$ws.Cells.Item(685, 6).Value = 45828.71103009259

# Row 734
$ws.Cells.Item(734, 3).Value = 24
$ws.Cells.Item(734, 4).Value = 45831.31672867793
$ws.Cells.Item(734, 5).Value = 24
$ws.Cells.Item(734, 6).Value = 45829.5324537037

# Row 772
$ws.Cells.Item(772, 3).Value = 468
$ws.Cells.Item(772, 4).Value = 45829.29943060575
$ws.Cells.Item(772, 5).Value = 468
$ws.Cells.Item(772, 6).Value = 45828.68149305556

# Row 839
$ws.Cells.Item(839, 3).Value = 19
$ws.Cells.Item(839, 4).Value = 45831.31672868106
$ws.Cells.Item(839, 5).Value = 19
$ws.Cells.Item(839, 6).Value = 45829.55277777778

# Row 889
$ws.Cells.Item(889, 3).Value = 24
$ws.Cells.Item(889, 4).Value = 45829.29943059556
$ws.Cells.Item(889, 5).Value = 24
$ws.Cells.Item(889, 6).Value = 45828.50344907407

# Row 929
$ws.Cells.Item(929, 3).Value = 49
$ws.Cells.Item(929, 4).Value = 45829.29943060546
$ws.Cells.Item(929, 5).Value = 49
$ws.Cells.Item(929, 6).Value = 45828.67482638889
$ws.Cells.Item(929, 7).Value = 0
$ws.Cells.Item(929, 8).Value = "Consistente"

# Row 963
$ws.Cells.Item(963, 3).Value = 1793
$ws.Cells.Item(963, 4).Value = 45829.29943059765
$ws.Cells.Item(963, 5).Value = 1793
$ws.Cells.Item(963, 6).Value = 45828.52178240741

# Row 1048
$ws.Cells.Item(1048, 3).Value = 26
$ws.Cells.Item(1048, 4).Value = 45831.31672868093
$ws.Cells.Item(1048, 5).Value = 26
$ws.Cells.Item(1048, 6).Value = 45829.55180555556

# Row 1092
$ws.Cells.Item(1092, 3).Value = 5
$ws.Cells.Item(1092, 4).Value = 45829.29943059825
$ws.Cells.Item(1092, 5).Value = 5
$ws.Cells.Item(1092, 6).Value = 45828.53510416667

# Row 1121
$ws.Cells.Item(1121, 3).Value = 28
$ws.Cells.Item(1121, 4).Value = 45829.29943059646
$ws.Cells.Item(1121, 5).Value = 28
$ws.Cells.Item(1121, 6).Value = 45828.52141203704

# Row 1135
$ws.Cells.Item(1135, 3).Value = 3
$ws.Cells.Item(1135, 4).Value = 45829.29943060112
$ws.Cells.Item(1135, 5).Value = 3
$ws.Cells.Item(1135, 6).Value = 45828.62064814815

# Row 1150
$ws.Cells.Item(1150, 3).Value = 27
$ws.Cells.Item(1150, 4).Value = 45829.29943060463
$ws.Cells.Item(1150, 5).Value = 27
$ws.Cells.Item(1150, 6).Value = 45828.66594907407

# Row 1157
$ws.Cells.Item(1157, 4).Value = 45831.31672868001

# Row 1176
$ws.Cells.Item(1176, 3).Value = 140
$ws.Cells.Item(1176, 4).Value = 45831.31672867721
$ws.Cells.Item(1176, 5).Value = 140
$ws.Cells.Item(1176, 6).Value = 45829.4890625

# Row 1177
$ws.Cells.Item(1177, 3).Value = 172
$ws.Cells.Item(1177, 4).Value = 45831.31672867704
$ws.Cells.Item(1177, 5).Value = 172
$ws.Cells.Item(1177, 6).Value = 45829.48858796297

# Row 1181
$ws.Cells.Item(1181, 3).Value = 47
$ws.Cells.Item(1181, 4).Value = 45831.31672868273
$ws.Cells.Item(1181, 5).Value = 47
$ws.Cells.Item(1181, 6).Value = 45829.58006944445

# Row 1183
$ws.Cells.Item(1183, 3).Value = 34
$ws.Cells.Item(1183, 4).Value = 45829.29943060681
$ws.Cells.Item(1183, 5).Value = 34
$ws.Cells.Item(1183, 6).Value = 45828.71103009259

# Row 1225
$ws.Cells.Item(1225, 3).Value = 11
$ws.Cells.Item(1225, 4).Value = 45831.31672868337
$ws.Cells.Item(1225, 5).Value = 11
$ws.Cells.Item(1225, 6).Value = 45829.58974537037

# Row 1226
$ws.Cells.Item(1226, 3).Value = 12
$ws.Cells.Item(1226, 4).Value = 45831.31672868355
$ws.Cells.Item(1226, 5).Value = 12
$ws.Cells.Item(1226, 6).Value = 45829.58974537037

# Row 1415
$ws.Cells.Item(1415, 3).Value = 100
$ws.Cells.Item(1415, 4).Value = 45831.31672867735
$ws.Cells.Item(1415, 5).Value = 100
$ws.Cells.Item(1415, 6).Value = 45829.48931712963

# Row 1555
$ws.Cells.Item(1555, 3).Value = 50
$ws.Cells.Item(1555, 4).Value = 45831.3167286775
$ws.Cells.Item(1555, 5).Value = 50
$ws.Cells.Item(1555, 6).Value = 45829.49069444444

# Row 1597
$ws.Cells.Item(1597, 3).Value = 5311
$ws.Cells.Item(1597, 4).Value = 45829.29943060262
$ws.Cells.Item(1597, 5).Value = 5311
$ws.Cells.Item(1597, 6).Value = 45828.63898148148

# Row 1650
$ws.Cells.Item(1650, 3).Value = -1
$ws.Cells.Item(1650, 4).Value = 45829.29943060434
$ws.Cells.Item(1650, 5).Value = -1
$ws.Cells.Item(1650, 6).Value = 45828.6547337963

# Row 1856
$ws.Cells.Item(1856, 3).Value = 11
$ws.Cells.Item(1856, 4).Value = 45829.29943059676
$ws.Cells.Item(1856, 5).Value = 11
$ws.Cells.Item(1856, 6).Value = 45828.52141203704

# Row 1876
$ws.Cells.Item(1876, 3).Value = 38
$ws.Cells.Item(1876, 4).Value = 45831.31672867807
$ws.Cells.Item(1876, 5).Value = 38
$ws.Cells.Item(1876, 6).Value = 45829.53827546296

# Row 1898
$ws.Cells.Item(1898, 3).Value = 0
$ws.Cells.Item(1898, 4).Value = 45829.29943059856
$ws.Cells.Item(1898, 5).Value = 0
$ws.Cells.Item(1898, 6).Value = 45828.54282407407

# Row 1941
$ws.Cells.Item(1941, 3).Value = 4
$ws.Cells.Item(1941, 4).Value = 45829.29943059885
$ws.Cells.Item(1941, 5).Value = 4
$ws.Cells.Item(1941, 6).Value = 45828.54484953704

# Row 2004
$ws.Cells.Item(2004, 3).Value = 4
$ws.Cells.Item(2004, 4).Value = 45829.29943060492
$ws.Cells.Item(2004, 5).Value = 4
$ws.Cells.Item(2004, 6).Value = 45828.66594907407

# Row 2099
$ws.Cells.Item(2099, 3).Value = 50
$ws.Cells.Item(2099, 4).Value = 45829.29943060518
$ws.Cells.Item(2099, 5).Value = 50
$ws.Cells.Item(2099, 6).Value = 45828.66594907407

# Row 2190
$ws.Cells.Item(2190, 3).Value = 9
$ws.Cells.Item(2190, 4).Value = 45831.31672867568
$ws.Cells.Item(2190, 5).Value = 9
$ws.Cells.Item(2190, 6).Value = 45829.47851851852

# Row 2191
$ws.Cells.Item(2191, 3).Value = 12
$ws.Cells.Item(2191, 4).Value = 45831.31672867687
$ws.Cells.Item(2191, 5).Value = 12
$ws.Cells.Item(2191, 6).Value = 45829.47966435185

# Row 2193
$ws.Cells.Item(2193, 3).Value = 11
$ws.Cells.Item(2193, 4).Value = 45831.31672867657
$ws.Cells.Item(2193, 5).Value = 11
$ws.Cells.Item(2193, 6).Value = 45829.47924768519

# Row 2194
$ws.Cells.Item(2194, 3).Value = 8
$ws.Cells.Item(2194, 4).Value = 45831.31672867673
$ws.Cells.Item(2194, 5).Value = 8
$ws.Cells.Item(2194, 6).Value = 45829.47944444444

# Row 2210
$ws.Cells.Item(2210, 3).Value = 11
$ws.Cells.Item(2210, 4).Value = 45831.31672867642
$ws.Cells.Item(2210, 5).Value = 11
$ws.Cells.Item(2210, 6).Value = 45829.47900462963

# Row 2264
$ws.Cells.Item(2264, 3).Value = 149
$ws.Cells.Item(2264, 4).Value = 45829.29943059482
$ws.Cells.Item(2264, 5).Value = 149
$ws.Cells.Item(2264, 6).Value = 45828.49071759259

# Row 2481
$ws.Cells.Item(2481, 3).Value = 398
$ws.Cells.Item(2481, 4).Value = 45829.29943060291
$ws.Cells.Item(2481, 5).Value = 398
$ws.Cells.Item(2481, 6).Value = 45828.63898148148

# Row 2482
$ws.Cells.Item(2482, 3).Value = 568
$ws.Cells.Item(2482, 4).Value = 45829.2994306032
$ws.Cells.Item(2482, 5).Value = 568
$ws.Cells.Item(2482, 6).Value = 45828.63898148148

# Row 2483
$ws.Cells.Item(2483, 3).Value = 373
$ws.Cells.Item(2483, 4).Value = 45829.29943060348
$ws.Cells.Item(2483, 5).Value = 373
$ws.Cells.Item(2483, 6).Value = 45828.63898148148

# Row 2484
$ws.Cells.Item(2484, 3).Value = 747
$ws.Cells.Item(2484, 4).Value = 45829.29943060376
$ws.Cells.Item(2484, 5).Value = 747
$ws.Cells.Item(2484, 6).Value = 45828.63898148148

# Row 2497
$ws.Cells.Item(2497, 4).Value = 45831.31672868244

# --- Append new rows at the end of the table ---
# Row 2579
$ws.Cells.Item(2579, 1).Value = 43804644
$ws.Cells.Item(2579, 2).Value = 1
$ws.Cells.Item(2579, 3).Value = 20
$ws.Cells.Item(2579, 4).Value = 45831.31672868061
$ws.Cells.Item(2579, 5).Value = 20
$ws.Cells.Item(2579, 6).Value = 45829.54894675926
$ws.Cells.Item(2579, 7).Value = 0
$ws.Cells.Item(2579, 8).Value = "Consistente"
$ws.Cells.Item(2579, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2579, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2580
$ws.Cells.Item(2580, 1).Value = 43804650
$ws.Cells.Item(2580, 2).Value = 1
$ws.Cells.Item(2580, 3).Value = 48
$ws.Cells.Item(2580, 4).Value = 45831.31672868065
$ws.Cells.Item(2580, 5).Value = 48
$ws.Cells.Item(2580, 6).Value = 45829.55087962963
$ws.Cells.Item(2580, 7).Value = 0
$ws.Cells.Item(2580, 8).Value = "Consistente"
$ws.Cells.Item(2580, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2580, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Applied all changes from diff"
